$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 8-10 (the "ECs" sending-cluster block no longer present)
$ws.Rows("8:10").Delete()

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Sema3c"
$ws.Range("C2").Value = "Nrp1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 39.09670133333334
$ws.Range("H2").Value = 117.290104
$ws.Range("I2").Value = 0.9758026125363394
$ws.Range("J2").Value = 0.9758026125363395
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 110.642708
$ws.Range("N2").Value = 331.928124
$ws.Range("O2").Value = 0.5476418925386564
$ws.Range("P2").Value = 0.5476418925386564
$ws.Range("Q2").Value = 4325.764909387211
$ws.Range("R2").Value = 38931.88418448489
$ws.Range("S2").Value = 0.5343903894735661
$ws.Range("T2").Value = 0.5343903894735662

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Sema3c"
$ws.Range("C3").Value = "Nrp1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 39.09670133333334
$ws.Range("H3").Value = 117.290104
$ws.Range("I3").Value = 0.9758026125363394
$ws.Range("J3").Value = 0.9758026125363395
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 63.66262833333334
$ws.Range("N3").Value = 190.987885
$ws.Range("O3").Value = 0.3151072754333865
$ws.Range("P3").Value = 0.3151072754333865
$ws.Range("Q3").Value = 2488.998766043338
$ws.Range("R3").Value = 22400.98889439004
$ws.Range("S3").Value = 0.3074825025971064
$ws.Range("T3").Value = 0.3074825025971065

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Sema3c"
$ws.Range("C4").Value = "Nrp1"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 39.09670133333334
$ws.Range("H4").Value = 117.290104
$ws.Range("I4").Value = 0.9758026125363394
$ws.Range("J4").Value = 0.9758026125363395
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 27.72944133333333
$ws.Range("N4").Value = 83.188324
$ws.Range("O4").Value = 0.1372508320279571
$ws.Range("P4").Value = 0.1372508320279571
$ws.Range("Q4").Value = 1084.129685949522
$ws.Range("R4").Value = 9757.167173545697
$ws.Range("S4").Value = 0.1339297204656668
$ws.Range("T4").Value = 0.1339297204656669

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Sema3c"
$ws.Range("C5").Value = "Nrp1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9694973333333333
$ws.Range("H5").Value = 2.908492
$ws.Range("I5").Value = 0.02419738746366056
$ws.Range("J5").Value = 0.02419738746366056
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 110.642708
$ws.Range("N5").Value = 331.928124
$ws.Range("O5").Value = 0.5476418925386564
$ws.Range("P5").Value = 0.5476418925386564
$ws.Range("Q5").Value = 107.2678103587786
$ws.Range("R5").Value = 965.4102932290078
$ws.Range("S5").Value = 0.01325150306509022
$ws.Range("T5").Value = 0.01325150306509022

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Sema3c"
$ws.Range("C6").Value = "Nrp1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9694973333333333
$ws.Range("H6").Value = 2.908492
$ws.Range("I6").Value = 0.02419738746366056
$ws.Range("J6").Value = 0.02419738746366056
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 63.66262833333334
$ws.Range("N6").Value = 190.987885
$ws.Range("O6").Value = 0.3151072754333865
$ws.Range("P6").Value = 0.3151072754333865
$ws.Range("Q6").Value = 61.72074840215778
$ws.Range("R6").Value = 555.48673561942
$ws.Range("S6").Value = 0.00762477283628006
$ws.Range("T6").Value = 0.00762477283628006

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Sema3c"
$ws.Range("C7").Value = "Nrp1"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9694973333333333
$ws.Range("H7").Value = 2.908492
$ws.Range("I7").Value = 0.02419738746366056
$ws.Range("J7").Value = 0.02419738746366056
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 27.72944133333333
$ws.Range("N7").Value = 83.188324
$ws.Range("O7").Value = 0.1372508320279571
$ws.Range("P7").Value = 0.1372508320279571
$ws.Range("Q7").Value = 26.88361942748977
$ws.Range("R7").Value = 241.952574847408
$ws.Range("S7").Value = 0.00332111156229027
$ws.Range("T7").Value = 0.00332111156229027
